$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "H9.2H"
$ws.Range("C3").Value = "H6.2H"
$ws.Range("C4").Value = "H5.2H"
$ws.Range("C5").Value = "H3.2H"
$ws.Range("C6").Value = "H2.2H"
$ws.Range("C7").Value = "A3.2C"
$ws.Range("C8").Value = "A5.2C"
$ws.Range("C9").Value = "A8.2C"
$ws.Range("C10").Value = "A6.2C"
$ws.Range("C11").Value = "A4.2C"
$ws.Range("C12").Value = "H8.2L"
$ws.Range("C13").Value = "H9.2L"
$ws.Range("C14").Value = "H3.2L"
$ws.Range("C15").Value = "H5.2L"
$ws.Range("C16").Value = "H2.2L"
$ws.Range("C17").Value = "A2.2H"
$ws.Range("C18").Value = "A1.2H"
$ws.Range("C19").Value = "A5.2H"
$ws.Range("C20").Value = "A8.2H"
$ws.Range("C21").Value = "A4.2H"
$ws.Range("C22").Value = "H10.2C"
$ws.Range("C23").Value = "H8.2C"
$ws.Range("C24").Value = "H2.2C"
$ws.Range("C25").Value = "H1.2C"
$ws.Range("C26").Value = "H3.2C"
$ws.Range("C27").Value = "A7.2L"
$ws.Range("C28").Value = "A2.2L"
$ws.Range("C29").Value = "A1.2L"
$ws.Range("C30").Value = "A5.2L"
$ws.Range("C31").Value = "A8.2L"
$ws.Range("C32").Value = "H5.2C"
$ws.Range("C33").Value = "H4.2C"
$ws.Range("C34").Value = "H9.2C"
$ws.Range("C35").Value = "H7.2C"
$ws.Range("C36").Value = "H6.2C"
$ws.Range("C37").Value = "A10.2L"
$ws.Range("C38").Value = "A9.2L"
$ws.Range("C39").Value = "A4.2L"
$ws.Range("C40").Value = "A3.2L"
$ws.Range("C41").Value = "A6.2L"
$ws.Range("C42").Value = "H4.2L"
$ws.Range("C43").Value = "H10.2L"
$ws.Range("C44").Value = "H1.2L"
$ws.Range("C45").Value = "H6.2L"
$ws.Range("C46").Value = "H7.2L"
$ws.Range("C47").Value = "A9.2H"
$ws.Range("C48").Value = "A7.2H"
$ws.Range("C49").Value = "A3.2H"
$ws.Range("C50").Value = "A6.2H"
$ws.Range("C51").Value = "A10.2H"
$ws.Range("C52").Value = "H7.2H"
$ws.Range("C53").Value = "H10.2H"
$ws.Range("C54").Value = "H4.2H"
$ws.Range("C55").Value = "H8.2H"
$ws.Range("C56").Value = "H1.2H"
$ws.Range("C57").Value = "A7.2C"
$ws.Range("C58").Value = "A9.2C"
$ws.Range("C59").Value = "A2.2C"
$ws.Range("C60").Value = "A10.2C"
$ws.Range("C61").Value = "A1.2C"

# Update the frozen-pane scroll position and active selection to match the authored view state.
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A2").Select()
$win.FreezePanes = $true
$ws.Range("C62").Select()
